$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "pi" in C1 and the value of pi in C2
$ws.Range("C1").Value = "pi"
$ws.Range("C2").Value = 3.14159265359

# Size column C to (best)fit its new contents, like Excel's AutoFit would
$ws.Columns.Item(3).ColumnWidth = 8.86

# Update the active selection to D2 as in the target workbook
$ws.Range("D2").Select() | Out-Null
